$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = "2020-12-25 00:00:00"
$ws.Range("N2").Value = "2019-12-31 00:00:00"

$ws.Range("O2").Value = 267978244.7
$ws.Range("P2").Value = 2366743.57
$ws.Range("Q2").Value = 125992468.46
$ws.Range("R2").Value = 305.7587276223
$ws.Range("S2").Value = 66829709.97
$ws.Range("T2").Value = 34.4590737332
$ws.Range("U2").Value = 36169359.51
$ws.Range("V2").Value = 6.1949635681
$ws.Range("W2").Value = 44062196.42
$ws.Range("X2").Value = 26582806.35
$ws.Range("Y2").Value = 21.4180695448
$ws.Range("Z2").Value = 1500296.7
$ws.Range("AA2").Value = -64.264334803
$ws.Range("AB2").Value = 223916048.28
$ws.Range("AC2").Value = 85.8515711889
$ws.Range("AD2").Value = 68.4135918099
$ws.Range("AE2").Value = 14.0384580679
$ws.Range("AF2").Value = 597.3124881277
$ws.Range("AG2").Value = 16.4424528078
